# major accuracy check update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the polyA isolation protocol kit name (shared string used in column G,
# rows 2:41) from "NEBNextPoly(A)E7490" to "NEBNextPoly(A)E7490L".
$ws.Range("G2:G41").Value = "NEBNextPoly(A)E7490L"

# Widen column G to fit the longer text.
$ws.Columns.Item(7).ColumnWidth = 24.03

# Reset the view: scroll back to the top-left and move the active
# selection to column G instead of column I.
[void]$ws.Range("G2:G41").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
